$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1950394588500564
$ws.Range("C2").Value = 0.5591882750845547
$ws.Range("J2").Value = 0.02029312288613303
$ws.Range("P2").Value = 0.1476888387824126
$ws.Range("S2").Value = 0.0777903043968433
$ws.Range("B3").Value = 0.007648183556405353
$ws.Range("C3").Value = 0.0401529636711281
$ws.Range("J3").Value = 0.03824091778202677
$ws.Range("P3").Value = 0.739961759082218
$ws.Range("S3").Value = 0.1739961759082218
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.7232142857142857
$ws.Range("S4").Value = 0.2053571428571428
$ws.Range("P5").Value = 0.7777777777777778
$ws.Range("S5").Value = 0.2222222222222222
$ws.Range("B6").Value = 0.06131386861313869
$ws.Range("D6").Value = 0.01021897810218978
$ws.Range("E6").Value = 0.00145985401459854
$ws.Range("F6").Value = 0.08613138686131387
$ws.Range("J6").Value = 0.1941605839416058
$ws.Range("O6").Value = 0.03211678832116788
$ws.Range("Q6").Value = 0.162043795620438
$ws.Range("R6").Value = 0.08321167883211679
$ws.Range("S6").Value = 0.3693430656934307
$ws.Range("B7").Value = 0.1099476439790576
$ws.Range("D7").Value = 0.02094240837696335
$ws.Range("F7").Value = 0.07504363001745201
$ws.Range("J7").Value = 0.1308900523560209
$ws.Range("O7").Value = 0.01745200698080279
$ws.Range("Q7").Value = 0.1972076788830716
$ws.Range("R7").Value = 0.06108202443280977
$ws.Range("S7").Value = 0.387434554973822
$ws.Range("B8").Value = 0.1069023569023569
$ws.Range("D8").Value = 0.01599326599326599
$ws.Range("E8").Value = 0.001683501683501683
$ws.Range("F8").Value = 0.0707070707070707
$ws.Range("J8").Value = 0.1144781144781145
$ws.Range("O8").Value = 0.01767676767676768
$ws.Range("Q8").Value = 0.1927609427609428
$ws.Range("R8").Value = 0.08585858585858586
$ws.Range("S8").Value = 0.3939393939393939
$ws.Range("B9").Value = 0.09611829944547134
$ws.Range("D9").Value = 0.01848428835489834
$ws.Range("F9").Value = 0.07578558225508318
$ws.Range("J9").Value = 0.08687615526802218
$ws.Range("O9").Value = 0.022181146025878
$ws.Range("Q9").Value = 0.1682070240295749
$ws.Range("R9").Value = 0.11090573012939
$ws.Range("S9").Value = 0.4214417744916821
$ws.Range("B10").Value = 0.1115281501340483
$ws.Range("D10").Value = 0.01769436997319035
$ws.Range("E10").Value = 0.00160857908847185
$ws.Range("F10").Value = 0.07024128686327077
$ws.Range("J10").Value = 0.1168900804289544
$ws.Range("O10").Value = 0.01715817694369973
$ws.Range("Q10").Value = 0.2235924932975871
$ws.Range("R10").Value = 0.08150134048257372
$ws.Range("S10").Value = 0.3597855227882037
$ws.Range("G11").Value = 0.1528795811518324
$ws.Range("J11").Value = 0.09214659685863874
$ws.Range("K11").Value = 0.2104712041884817
$ws.Range("L11").Value = 0.5319371727748691
$ws.Range("S11").Value = 0.01256544502617801
$ws.Range("G12").Value = 0.6836158192090396
$ws.Range("J12").Value = 0.2184557438794727
$ws.Range("K12").Value = 0.02448210922787194
$ws.Range("L12").Value = 0.04143126177024482
$ws.Range("S12").Value = 0.032015065913371
$ws.Range("G13").Value = 0.6611570247933884
$ws.Range("J13").Value = 0.3140495867768595
$ws.Range("S13").Value = 0.02479338842975207
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.02009273570324575
$ws.Range("H15").Value = 0.1437403400309119
$ws.Range("I15").Value = 0.05100463678516229
$ws.Range("J15").Value = 0.3446676970633694
$ws.Range("K15").Value = 0.1004636785162288
$ws.Range("M15").Value = 0.0170015455950541
$ws.Range("O15").Value = 0.07727975270479134
$ws.Range("S15").Value = 0.2457496136012365
$ws.Range("F16").Value = 0.0273037542662116
$ws.Range("H16").Value = 0.1467576791808874
$ws.Range("I16").Value = 0.06996587030716724
$ws.Range("J16").Value = 0.4163822525597269
$ws.Range("K16").Value = 0.1160409556313993
$ws.Range("M16").Value = 0.0204778156996587
$ws.Range("N16").Value = 0.001706484641638225
$ws.Range("O16").Value = 0.06143344709897611
$ws.Range("S16").Value = 0.1399317406143345
$ws.Range("F17").Value = 0.01297764960346071
$ws.Range("H17").Value = 0.1751982696467195
$ws.Range("I17").Value = 0.09805335255948089
$ws.Range("J17").Value = 0.4232155731795241
$ws.Range("K17").Value = 0.08940158615717375
$ws.Range("M17").Value = 0.0144196106705119
$ws.Range("N17").Value = 0.002162941600576784
$ws.Range("O17").Value = 0.05695746214852199
$ws.Range("S17").Value = 0.1276135544340303
$ws.Range("F18").Value = 0.01971326164874552
$ws.Range("H18").Value = 0.1756272401433692
$ws.Range("I18").Value = 0.08781362007168458
$ws.Range("J18").Value = 0.4103942652329749
$ws.Range("K18").Value = 0.0985663082437276
$ws.Range("M18").Value = 0.02150537634408602
$ws.Range("O18").Value = 0.06630824372759857
$ws.Range("S18").Value = 0.1200716845878136
$ws.Range("F19").Value = 0.01402254605444047
$ws.Range("H19").Value = 0.1875171844927138
$ws.Range("I19").Value = 0.076711575474292
$ws.Range("J19").Value = 0.380808358537256
$ws.Range("K19").Value = 0.1146549353863074
$ws.Range("M19").Value = 0.0195215837228485
$ws.Range("N19").Value = 0.001374759417102007
$ws.Range("O19").Value = 0.07011273027220237
$ws.Range("S19").Value = 0.1352763266428375
